# Week 13 logging update for Chargers 2021 Team Data workbook.
# Appends this week's per-play logs to the running season logs on YDS/ST,
# and rolls the season-to-date totals forward on OFF/DEF/TURNS/PEN.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# YDS — running play-by-play yardage logs (Home/Road x R/P columns)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = $ws.Range("B2").Value() + " 2 0 -1 6 7 3 9 5 0 2 5 4 1 15 1 5 6 1 3 -4 7 2"
$ws.Range("C2").Value = $ws.Range("C2").Value() + " 0 2 2 7 1 -5 6 8 6 2 6 17 2 5 3 4 0 10 3 3 8 -3 5 10"
$ws.Range("B3").Value = $ws.Range("B3").Value() + " -1 10 -1 41 5 10 47 6 7 44 23 7 4 7 12 7 12 9 7 6 -3 1 18 33 2"
$ws.Range("C3").Value = $ws.Range("C3").Value() + " 5 12 21 3 9 7 32 16 29 5 7 0 5 9 20 12 6 9 11 12 16 14 25 14 10"

# ---------------------------------------------------------------------
# OFF — season offensive totals (row2 = Home, row3 = Road)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("OFF")

$ws.Range("C2").Value = 126
$ws.Range("D2").Value = 10
$ws.Range("E2").Value = 6
$ws.Range("F2").Value = 32
$ws.Range("G2").Value = 37
$ws.Range("I2").Value = 6
$ws.Range("J2").Value = 11
$ws.Range("N2").Value = 12
$ws.Range("O2").Value = 21
$ws.Range("P2").Value = 13

$ws.Range("B3").Value = 8
$ws.Range("C3").Value = 155
$ws.Range("E3").Value = 30
$ws.Range("F3").Value = 81
$ws.Range("G3").Value = 36
$ws.Range("H3").Value = 24
$ws.Range("I3").Value = 44
$ws.Range("J3").Value = 48
$ws.Range("L3").Value = 240
$ws.Range("M3").Value = 165
$ws.Range("Q3").Value = 380

# ---------------------------------------------------------------------
# DEF — season defensive totals (row2 = Home, row3 = Road)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("DEF")

$ws.Range("C2").Value = 174
$ws.Range("D2").Value = 11
$ws.Range("F2").Value = 45
$ws.Range("G2").Value = 56
$ws.Range("N2").Value = 14
$ws.Range("O2").Value = 14

$ws.Range("B3").Value = 8
$ws.Range("C3").Value = 118
$ws.Range("E3").Value = 28
$ws.Range("F3").Value = 64
$ws.Range("G3").Value = 36
$ws.Range("H3").Value = 17
$ws.Range("I3").Value = 36
$ws.Range("J3").Value = 41
$ws.Range("L3").Value = 174
$ws.Range("M3").Value = 111
$ws.Range("Q3").Value = 382

# ---------------------------------------------------------------------
# ST — special teams totals + running logs
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ST")

$ws.Range("B2").Value = 64
$ws.Range("D2").Value = 38
$ws.Range("F2").Value = 210
$ws.Range("G2").Value = 197
$ws.Range("H2").Value = 7
$ws.Range("I2").Value = 5
$ws.Range("L2").Value = 74
$ws.Range("M2").Value = 61

$ws.Range("B3").Value = 35

$ws.Range("D3").Value = $ws.Range("D3").Value() + " 39 49 45 36"
$ws.Range("B4").Value = $ws.Range("B4").Value() + " 66 63 62"
$ws.Range("D4").Value = $ws.Range("D4").Value() + " 0 5 0 0"
$ws.Range("B5").Value = $ws.Range("B5").Value() + " 17 10 25"
$ws.Range("D5").Value = "0 0"
$ws.Range("B6").Value = $ws.Range("B6").Value() + " 48 24 27"

# ---------------------------------------------------------------------
# TURNS — season turnover totals (row3 = Road)
# ---------------------------------------------------------------------
$turns = $wb.Worksheets.Item("TURNS")

$turns.Range("B3").Value = 5
$turns.Range("C3").Value = 7
$turns.Range("D3").Value = 4
$turns.Range("E3").Value = 8

# ---------------------------------------------------------------------
# PEN — season penalty totals
# ---------------------------------------------------------------------
$pen = $wb.Worksheets.Item("PEN")

$pen.Range("B2").Value = 13
$pen.Range("D2").Value = 8
$pen.Range("B3").Value = 18
$pen.Range("B4").Value = 4

# ---------------------------------------------------------------------
# Leave the ST sheet active with D6 selected, matching the last thing
# touched during this week's logging session.
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("D6").Select()
